$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "Meg Davis" (WCP 4.134) in sorted position, between
# "Meg Davidson" (row 16) and "Ashwini Deo" (row 17 before the insert).
$ws.Rows.Item(17).Insert()
$ws.Cells.Item(17, 1).Value = "Meg"
$ws.Cells.Item(17, 2).Value = "Davis"
$ws.Cells.Item(17, 3).Value = "WCP 4.134"
$ws.Cells.Item(17, 4).Formula = "=CONCATENATE(A17, "" "", B17, ""   "", C17)"

# Update Baorian Nuchged's office/label to include the second room (now row 43
# after the insertion above shifted everything down by one).
$ws.Cells.Item(43, 3).Value = "4.110/W8"

# Fix a misspelled first name: "Tiago" -> "Thiago" (Cardoso Aguiar, row 11)
$ws.Cells.Item(11, 1).Value = "Thiago"

# Re-apply the existing last-name sort so the sheet's remembered sort range
# grows to include the newly inserted row (A3:D67, keyed on column B).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B3:B67"))
$sortObj.SetRange($ws.Range("A3:D67"))
$sortObj.Header = -4142
$sortObj.Apply()

# Restore the default view: no frozen/scrolled top-left cell, selection at A11.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A11").Select()
